# Apply the latest cryptocurrency market-data refresh (price + 1h volume change).
# GitHub Actions scheduled job: pulls fresh rankings and rewrites the changed cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.117.39"
$ws.Range("E2").Value = "  +4.79%  "
# Row 3
$ws.Range("D3").Value = "2.704.09"
$ws.Range("E3").Value = "  +3.93%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.60"
$ws.Range("E5").Value = "  +0.27%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.24"
$ws.Range("E6").Value = "  +4.34%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.35%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +1.43%  "
# Row 9
$ws.Range("D9").Value = "2.734.76"
$ws.Range("E9").Value = "  +4.90%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.74"
$ws.Range("E10").Value = "  +2.65%  "
# Row 11
$ws.Range("E11").Value = "  +6.86%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  +4.24%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.158"
$ws.Range("E13").Value = "  +1.41%  "
# Row 14
$ws.Range("D14").Value = "3.186.63"
$ws.Range("E14").Value = "  +4.07%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.69"
$ws.Range("E15").Value = "  +9.23%  "
# Row 16
$ws.Range("D16").Value = "63.010.58"
$ws.Range("E16").Value = "  +4.61%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000150"
$ws.Range("E17").Value = "  +7.00%  "
# Row 18
$ws.Range("D18").Value = "2.720.70"
$ws.Range("E18").Value = "  +4.47%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.94"
$ws.Range("E19").Value = "  +5.09%  "
# Row 20
$ws.Range("E20").Value = "  +5.27%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "362.29"
$ws.Range("E21").Value = "  +4.79%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.00"
$ws.Range("E22").Value = "  +1.47%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.40%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.530"
$ws.Range("E24").Value = "  -0.65%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.42"
$ws.Range("E25").Value = "  +2.80%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  +3.83%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.65"
$ws.Range("E27").Value = "  +8.12%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  -0.30%  "
# Row 29
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0854"
$ws.Range("E29").Value = "  +6.99%  "
# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.02"
$ws.Range("E30").Value = "  +5.63%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.06"
$ws.Range("E31").Value = "  +10.46%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.59"
$ws.Range("E32").Value = "  +1.53%  "
# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  +22.38%  "
# Row 34
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.19%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.52"
$ws.Range("E35").Value = "  +5.58%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.76"
$ws.Range("E36").Value = "  +11.91%  "
# Row 37
$ws.Range("E37").Value = "  +8.14%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.81"
$ws.Range("E38").Value = "  +10.38%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +18.89%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "349.49"
$ws.Range("E40").Value = "  +11.61%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.26"
$ws.Range("E41").Value = "  +9.74%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.22"
$ws.Range("E42").Value = "  +2.78%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.64"
$ws.Range("E43").Value = "  +13.60%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.58"
$ws.Range("E44").Value = "  +8.37%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0592"
$ws.Range("E45").Value = "  +7.82%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.64"
$ws.Range("E46").Value = "  +8.86%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0260"
$ws.Range("E47").Value = "  +6.86%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.55"
$ws.Range("E48").Value = "  +2.25%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.640"
$ws.Range("E49").Value = "  +5.60%  "
# Row 50
$ws.Range("E50").Value = "  +1.42%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.996"
$ws.Range("E51").Value = "  -0.37%  "
